$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder the block-type labels ---
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# --- Data rows (2-7): updated 0/1 indicator values ---
$data = @(
    @(0,0,0,0,1,0),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,1,0,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
